$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers for the two new columns (I and J) ---
# Copy the format from the existing header cell H1 (bold font, border,
# centered alignment) onto the new header cells so they reuse the same
# style rather than Excel creating a brand-new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data for I2:J73 ---
$data = New-Object 'object[,]' 72,2
$data[0,0] = 7; $data[0,1] = 7
$data[1,0] = 6; $data[1,1] = 6
$data[2,0] = 7; $data[2,1] = 7
$data[3,0] = 10; $data[3,1] = 11
$data[4,0] = 6; $data[4,1] = 7
$data[5,0] = 8; $data[5,1] = 8
$data[6,0] = 7; $data[6,1] = 7
$data[7,0] = 6; $data[7,1] = 7
$data[8,0] = 7; $data[8,1] = 7
$data[9,0] = 7; $data[9,1] = 7
$data[10,0] = 8; $data[10,1] = 8
$data[11,0] = 8; $data[11,1] = 8
$data[12,0] = 9; $data[12,1] = 9
$data[13,0] = 7; $data[13,1] = 8
$data[14,0] = 7; $data[14,1] = 7
$data[15,0] = 5; $data[15,1] = 6
$data[16,0] = 7; $data[16,1] = 7
$data[17,0] = 7; $data[17,1] = 8
$data[18,0] = 7; $data[18,1] = 7
$data[19,0] = 7; $data[19,1] = 8
$data[20,0] = 7; $data[20,1] = 7
$data[21,0] = 4; $data[21,1] = 5
$data[22,0] = 6; $data[22,1] = 6
$data[23,0] = 9; $data[23,1] = 9
$data[24,0] = 6; $data[24,1] = 7
$data[25,0] = 9; $data[25,1] = 9
$data[26,0] = 9; $data[26,1] = 9
$data[27,0] = 9; $data[27,1] = 9
$data[28,0] = 8; $data[28,1] = 8
$data[29,0] = 9; $data[29,1] = 9
$data[30,0] = 9; $data[30,1] = 9
$data[31,0] = 9; $data[31,1] = 9
$data[32,0] = 8; $data[32,1] = 8
$data[33,0] = 8; $data[33,1] = 9
$data[34,0] = 9; $data[34,1] = 9
$data[35,0] = 9; $data[35,1] = 9
$data[36,0] = 9; $data[36,1] = 9
$data[37,0] = 9; $data[37,1] = 9
$data[38,0] = 8; $data[38,1] = 8
$data[39,0] = 9; $data[39,1] = 9
$data[40,0] = 8; $data[40,1] = 8
$data[41,0] = 9; $data[41,1] = 9
$data[42,0] = 9; $data[42,1] = 9
$data[43,0] = 9; $data[43,1] = 9
$data[44,0] = 9; $data[44,1] = 9
$data[45,0] = 9; $data[45,1] = 10
$data[46,0] = 9; $data[46,1] = 9
$data[47,0] = 9; $data[47,1] = 9
$data[48,0] = 8; $data[48,1] = 8
$data[49,0] = 9; $data[49,1] = 9
$data[50,0] = 9; $data[50,1] = 9
$data[51,0] = 9; $data[51,1] = 9
$data[52,0] = 9; $data[52,1] = 9
$data[53,0] = 9; $data[53,1] = 9
$data[54,0] = 8; $data[54,1] = 9
$data[55,0] = 9; $data[55,1] = 9
$data[56,0] = 9; $data[56,1] = 9
$data[57,0] = 9; $data[57,1] = 9
$data[58,0] = 9; $data[58,1] = 9
$data[59,0] = 10; $data[59,1] = 10
$data[60,0] = 9; $data[60,1] = 9
$data[61,0] = 9; $data[61,1] = 9
$data[62,0] = 9; $data[62,1] = 9
$data[63,0] = 5; $data[63,1] = 5
$data[64,0] = 4; $data[64,1] = 4
$data[65,0] = 6; $data[65,1] = 6
$data[66,0] = 4; $data[66,1] = 4
$data[67,0] = 4; $data[67,1] = 4
$data[68,0] = 8; $data[68,1] = 8
$data[69,0] = 4; $data[69,1] = 4
$data[70,0] = 4; $data[70,1] = 4
$data[71,0] = 3; $data[71,1] = 3

$ws.Range("I2:J73").Value = $data

$wb.Save()
